$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PESSOAL")

# Group 1: rows 7-31 -> ID 123956789, "Liliane Santoé", date serial 31569
$ws.Range("A7:A31").Value = 123956789
$ws.Range("B7:B31").Value = "Liliane Santoé"
$ws.Range("C7:C31").Value = 31569

# Group 2: rows 32-37 -> ID 122456789, "Sueli Santos '4", date serial 29871
$ws.Range("A32:A37").Value = 122456789
$ws.Range("B32:B37").Value = "Sueli Santos '4"
$ws.Range("C32:C37").Value = 29871

# Group 3: row 38 -> ID 123956789, "Liliane Santoé", date serial 31569
$ws.Range("A38").Value = 123956789
$ws.Range("B38").Value = "Liliane Santoé"
$ws.Range("C38").Value = 31569

# Match the existing date-style (numFmt 14) used by C4:C6 onto the new rows
# by copying the format from an existing styled cell rather than assigning
# a NumberFormat directly (which would create a brand new style entry).
$ws.Range("C6").Copy()
$ws.Range("C7:C38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A2").Select()
